$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price rows for "Níspero" get reshuffled: the pair of rows
# currently on 2022-12-16 (serial 44911) swap places with the pair
# currently on 2022-12-20 (serial 44915), keeping the internal order of
# each pair (Primera/Segunda stays together, Especial/Primera stays
# together).
#
# Effect: row 2 <-> row 4, and row 3 <-> row 5, for columns D..S
# (the surrounding columns A,B,C,E..K,Q,T are identical in every row so
# swapping them is a no-op, but we include the full D:S span to mirror
# the diff precisely).

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"

    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $v4 = $ws.Range($addr4).Value2
    $v5 = $ws.Range($addr5).Value2

    $ws.Range($addr2).Value2 = $v4
    $ws.Range($addr4).Value2 = $v2

    $ws.Range($addr3).Value2 = $v5
    $ws.Range($addr5).Value2 = $v3
}
